$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1289.0549
$ws.Range("I15").Value = 1289.0549
$ws.Range("K15").Value = 3867.1647
$ws.Range("M15").Value = -3698.1647
$ws.Range("H40").Value = 8554.909
$ws.Range("I40").Value = 7469.25
$ws.Range("J40").Value = 11450
$ws.Range("K40").Value = 7469.25
$ws.Range("L40").Value = 11450
$ws.Range("M40").Value = -7294.25
$ws.Range("N40").Value = -11800
$ws.Range("H41").Value = 35718668
$ws.Range("J41").Value = 125013816
$ws.Range("L41").Value = 125013816
$ws.Range("N41").Value = -125014696
$ws.Range("H51").Value = 3499.5
$ws.Range("I51").Value = 1999
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 1999
$ws.Range("L51").Value = 5000
$ws.Range("M51").Value = -1515
$ws.Range("N51").Value = -5968
$ws.Range("H62").Value = 17863416
$ws.Range("I62").Value = 31254976
$ws.Range("K62").Value = 31254976
$ws.Range("M62").Value = -31254352
$ws.Range("H64").Value = 10333.333
$ws.Range("J64").Value = 11400
$ws.Range("L64").Value = 11400
$ws.Range("N64").Value = -11896
$ws.Range("H65").Value = 17863416
$ws.Range("I65").Value = 31254976
$ws.Range("K65").Value = 156274880
$ws.Range("M65").Value = -156271760
$ws.Range("H67").Value = 10333.333
$ws.Range("J67").Value = 11400
$ws.Range("L67").Value = 11400
$ws.Range("N67").Value = -13116
$ws.Range("H69").Value = 17372.5
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 17372.5
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H74").Value = 17831.666
$ws.Range("I74").Value = 3500
$ws.Range("K74").Value = 3500
$ws.Range("M74").Value = -2564
$ws.Range("H76").Value = 45460056
$ws.Range("I76").Value = 5110.4614
$ws.Range("J76").Value = 111117200
$ws.Range("K76").Value = 5110.4614
$ws.Range("L76").Value = 111117200
$ws.Range("M76").Value = -4795.4614
$ws.Range("N76").Value = -111117830
$ws.Range("H77").Value = 17831.666
$ws.Range("I77").Value = 3500
$ws.Range("K77").Value = 17500
$ws.Range("M77").Value = -12820
$ws.Range("H79").Value = 45460056
$ws.Range("I79").Value = 5110.4614
$ws.Range("J79").Value = 111117200
$ws.Range("K79").Value = 5110.4614
$ws.Range("L79").Value = 111117200
$ws.Range("M79").Value = -4018.4614
$ws.Range("N79").Value = -111119384
$ws.Range("H93").Value = 54300.5
$ws.Range("I93").Value = 78000
$ws.Range("K93").Value = 78000
$ws.Range("M93").Value = -75504
$ws.Range("H94").Value = 1801
$ws.Range("I94").Value = 1705.3334
$ws.Range("J94").Value = 2375
$ws.Range("K94").Value = 1705.3334
$ws.Range("L94").Value = 2375
$ws.Range("M94").Value = -1254.3334
$ws.Range("N94").Value = -3277
$ws.Range("H100").Value = 3740.2222
$ws.Range("I100").Value = 3023.3572
$ws.Range("J100").Value = 6249.25
$ws.Range("K100").Value = 3023.3572
$ws.Range("L100").Value = 6249.25
$ws.Range("M100").Value = -2482.3572
$ws.Range("N100").Value = -7331.25
$ws.Range("H101").Value = 726
$ws.Range("I101").Value = 827.4
$ws.Range("J101").Value = 641.5
$ws.Range("K101").Value = 2482.2
$ws.Range("L101").Value = 1924.5
$ws.Range("M101").Value = -860.1999999999998
$ws.Range("N101").Value = -5168.5
$ws.Range("H125").Value = 10104464
$ws.Range("I125").Value = 848
$ws.Range("J125").Value = 13893319
$ws.Range("K125").Value = 7632
$ws.Range("L125").Value = 125039871
$ws.Range("M125").Value = -5172
$ws.Range("N125").Value = -125044791
$ws.Range("H131").Value = 4064.8096
$ws.Range("I131").Value = 2768.2856
$ws.Range("J131").Value = 6657.857
$ws.Range("K131").Value = 8304.856800000001
$ws.Range("L131").Value = 19973.571
$ws.Range("M131").Value = -3264.856800000001
$ws.Range("N131").Value = -30053.571
$ws.Range("H132").Value = 3254.7693
$ws.Range("I132").Value = 2853.3044
$ws.Range("K132").Value = 8559.913199999999
$ws.Range("M132").Value = -6029.913199999999
$ws.Range("H137").Value = 3575.3235
$ws.Range("I137").Value = 2740.2917
$ws.Range("J137").Value = 5579.4
$ws.Range("K137").Value = 8220.875100000001
$ws.Range("L137").Value = 16738.2
$ws.Range("M137").Value = -5670.875100000001
$ws.Range("N137").Value = -21838.2
$ws.Range("H138").Value = 6043.6587
$ws.Range("I138").Value = 2744.4783
$ws.Range("J138").Value = 7267.5483
$ws.Range("K138").Value = 8233.4349
$ws.Range("L138").Value = 21802.6449
$ws.Range("M138").Value = -3093.4349
$ws.Range("N138").Value = -32082.6449
$ws.Range("H141").Value = 6279.96
$ws.Range("I141").Value = 6265.174
$ws.Range("K141").Value = 18795.522
$ws.Range("M141").Value = -13615.522
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 53729.895
$ws.Range("I2").Value = 59704
$ws.Range("K2").Value = 59704
$ws.Range("M2").Value = -59591
$ws.Range("H17").Value = 9999
$ws.Range("J17").Value = 9999
$ws.Range("L17").Value = 9999
$ws.Range("N17").Value = -10345
$ws.Range("H32").Value = 2883.47
$ws.Range("I32").Value = 2389.0312
$ws.Range("J32").Value = 14750
$ws.Range("K32").Value = 2389.0312
$ws.Range("L32").Value = 14750
$ws.Range("M32").Value = -2102.0312
$ws.Range("N32").Value = -15324
$ws.Range("H45").Value = 2862.4
$ws.Range("I45").Value = 2002.75
$ws.Range("J45").Value = 3435.5
$ws.Range("K45").Value = 2002.75
$ws.Range("L45").Value = 3435.5
$ws.Range("M45").Value = -1625.75
$ws.Range("N45").Value = -4189.5
$ws.Range("H46").Value = 3997.5
$ws.Range("I46").Value = 4000
$ws.Range("J46").Value = 3996.6667
$ws.Range("K46").Value = 4000
$ws.Range("L46").Value = 3996.6667
$ws.Range("M46").Value = -3681
$ws.Range("N46").Value = -4634.6667
$ws.Range("H61").Value = 4364.16
$ws.Range("I61").Value = 4364.16
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4364.16
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4152.16
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 2276.7778
$ws.Range("I74").Value = 2227.543
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 2227.543
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -1353.543
$ws.Range("N74").Value = -5748
$ws.Range("H77").Value = 2276.7778
$ws.Range("I77").Value = 2227.543
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 11137.715
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -6769.715
$ws.Range("N77").Value = -28736
$ws.Range("H97").Value = 2642.5715
$ws.Range("I97").Value = 2642.5715
$ws.Range("K97").Value = 2642.5715
$ws.Range("M97").Value = -2146.5715
$ws.Range("H110").Value = 202136.64
$ws.Range("I110").Value = 265374.62
$ws.Range("J110").Value = 1883
$ws.Range("K110").Value = 265374.62
$ws.Range("L110").Value = 1883
$ws.Range("M110").Value = -263329.62
$ws.Range("N110").Value = -5973
$ws.Range("H113").Value = 100000
$ws.Range("J113").Value = 100000
$ws.Range("L113").Value = 100000
$ws.Range("N113").Value = -108678
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 53729.895
$ws.Range("I116").Value = 59704
$ws.Range("K116").Value = 59704
$ws.Range("M116").Value = -57410
$ws.Range("H118").Value = 100000
$ws.Range("J118").Value = 100000
$ws.Range("L118").Value = 100000
$ws.Range("N118").Value = -103314
$ws.Range("H119").Value = 149250
$ws.Range("J119").Value = 149250
$ws.Range("L119").Value = 149250
$ws.Range("N119").Value = -158926
$ws.Range("H120").Value = 65000
$ws.Range("J120").Value = 65000
$ws.Range("L120").Value = 65000
$ws.Range("N120").Value = -74676
$ws.Range("H122").Value = 4221.0645
$ws.Range("I122").Value = 3023.7856
$ws.Range("J122").Value = 5207.0586
$ws.Range("K122").Value = 9071.356800000001
$ws.Range("L122").Value = 15621.1758
$ws.Range("M122").Value = -6621.356800000001
$ws.Range("N122").Value = -20521.1758
$ws.Range("H128").Value = 65000
$ws.Range("J128").Value = 65000
$ws.Range("L128").Value = 65000
$ws.Range("N128").Value = -74960
$ws.Range("H132").Value = 4108.2856
$ws.Range("I132").Value = 3347.8635
$ws.Range("K132").Value = 10043.5905
$ws.Range("M132").Value = -7513.5905
$ws.Range("H136").Value = 4364.16
$ws.Range("I136").Value = 4364.16
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 13092.48
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -10542.48
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 53729.895
$ws.Range("I3").Value = 59704
$ws.Range("K3").Value = 59704
$ws.Range("M3").Value = -59590
$ws.Range("H94").Value = 1310.75
$ws.Range("I94").Value = 1391.1765
$ws.Range("J94").Value = 855
$ws.Range("K94").Value = 1391.1765
$ws.Range("L94").Value = 855
$ws.Range("M94").Value = -940.1765
$ws.Range("N94").Value = -1757
$ws.Range("H105").Value = 169153
$ws.Range("I105").Value = 202183.8
$ws.Range("K105").Value = 202183.8
$ws.Range("M105").Value = -200436.8
$ws.Range("H107").Value = 1116046.5
$ws.Range("I107").Value = 4948.6
$ws.Range("J107").Value = 2504918.8
$ws.Range("K107").Value = 4948.6
$ws.Range("L107").Value = 2504918.8
$ws.Range("M107").Value = -3028.6
$ws.Range("N107").Value = -2508758.8
$ws.Range("H134").Value = 25233.283
$ws.Range("I134").Value = 3534.0715
$ws.Range("J134").Value = 253075
$ws.Range("K134").Value = 10602.2145
$ws.Range("L134").Value = 759225
$ws.Range("M134").Value = -8067.2145
$ws.Range("N134").Value = -764295
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 3760
$ws.Range("I13").Value = 3000
$ws.Range("K13").Value = 3000
$ws.Range("M13").Value = -2861
$ws.Range("H22").Value = 835.3333
$ws.Range("I22").Value = 502.44446
$ws.Range("K22").Value = 502.44446
$ws.Range("M22").Value = -152.44446
$ws.Range("H26").Value = 8008.1665
$ws.Range("J26").Value = 8008.1665
$ws.Range("L26").Value = 8008.1665
$ws.Range("N26").Value = -8582.166499999999
$ws.Range("H58").Value = 2957.2778
$ws.Range("I58").Value = 3050.9167
$ws.Range("J58").Value = 2770
$ws.Range("K58").Value = 3050.9167
$ws.Range("L58").Value = 2770
$ws.Range("M58").Value = -2847.9167
$ws.Range("N58").Value = -3176
$ws.Range("H81").Value = 16000
$ws.Range("I81").Value = 16000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 16000
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("M81").Value = -15002
$ws.Range("H84").Value = 16000
$ws.Range("I84").Value = 16000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 48000
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("M84").Value = -43008
$ws.Range("H94").Value = 1739
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 1923.75
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 1923.75
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -2825.75
$ws.Range("H99").Value = 8343.5
$ws.Range("I99").Value = 9116.666999999999
$ws.Range("J99").Value = 7879.6
$ws.Range("K99").Value = 9116.666999999999
$ws.Range("L99").Value = 7879.6
$ws.Range("M99").Value = -7618.666999999999
$ws.Range("N99").Value = -10875.6
$ws.Range("H105").Value = 863.9231
$ws.Range("I105").Value = 863.9231
$ws.Range("K105").Value = 863.9231
$ws.Range("M105").Value = 883.0769
$ws.Range("H107").Value = 323.26666
$ws.Range("I107").Value = 266.45456
$ws.Range("J107").Value = 479.5
$ws.Range("K107").Value = 266.45456
$ws.Range("L107").Value = 479.5
$ws.Range("M107").Value = 1653.54544
$ws.Range("N107").Value = -4319.5
$ws.Range("H114").Value = 74998
$ws.Range("J114").Value = 74998
$ws.Range("L114").Value = 74998
$ws.Range("N114").Value = -83676
$ws.Range("H126").Value = 8343.5
$ws.Range("I126").Value = 9116.666999999999
$ws.Range("J126").Value = 7879.6
$ws.Range("K126").Value = 27350.001
$ws.Range("L126").Value = 23638.8
$ws.Range("M126").Value = -24880.001
$ws.Range("N126").Value = -28578.8
$ws.Range("H134").Value = 502746.4
$ws.Range("I134").Value = 3157.5293
$ws.Range("K134").Value = 9472.5879
$ws.Range("M134").Value = -6937.5879
$ws.Range("H135").Value = 49994.184
$ws.Range("J135").Value = 49994.184
$ws.Range("L135").Value = 49994.184
$ws.Range("N135").Value = -60134.184
$ws.Range("H136").Value = 2957.2778
$ws.Range("I136").Value = 3050.9167
$ws.Range("J136").Value = 2770
$ws.Range("K136").Value = 9152.750100000001
$ws.Range("L136").Value = 8310
$ws.Range("M136").Value = -6602.750100000001
$ws.Range("N136").Value = -13410
$ws.Range("H139").Value = 99415.336
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 99415.336
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 99415.336
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -109695.336
$ws.Range("H141").Value = 197150.7
$ws.Range("J141").Value = 195688.83
$ws.Range("L141").Value = 195688.83
$ws.Range("N141").Value = -206048.83
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1080551.2
$ws.Range("I5").Value = 89500.89
$ws.Range("K5").Value = 268502.67
$ws.Range("M5").Value = -268390.67
$ws.Range("H33").Value = 2377412
$ws.Range("I33").Value = 7407420.5
$ws.Range("K33").Value = 44444523
$ws.Range("M33").Value = -44444240
$ws.Range("H105").Value = 34290.145
$ws.Range("J105").Value = 35005
$ws.Range("L105").Value = 105015
$ws.Range("N105").Value = -110257
$ws.Range("H107").Value = 108939.69
$ws.Range("J107").Value = 206020.1
$ws.Range("L107").Value = 618060.3
$ws.Range("N107").Value = -621900.3
$ws.Range("H109").Value = 38747.586
$ws.Range("I109").Value = 2690
$ws.Range("J109").Value = 50220.453
$ws.Range("K109").Value = 8070
$ws.Range("L109").Value = 150661.359
$ws.Range("M109").Value = -7030
$ws.Range("N109").Value = -152741.359
$ws.Range("H113").Value = 1426216
$ws.Range("J113").Value = 1808.0869
$ws.Range("L113").Value = 5424.2607
$ws.Range("N113").Value = -9764.260699999999
$ws.Range("H117").Value = 1372.2858
$ws.Range("I117").Value = 261.25
$ws.Range("K117").Value = 783.75
$ws.Range("M117").Value = 2658.25
$ws.Range("H121").Value = 34334864
$ws.Range("I121").Value = 1965.3334
$ws.Range("J121").Value = 49048964
$ws.Range("K121").Value = 5896.0002
$ws.Range("L121").Value = 147146892
$ws.Range("M121").Value = -4586.0002
$ws.Range("N121").Value = -147149512
$ws.Range("H127").Value = 1442.6364
$ws.Range("J127").Value = 1442.6364
$ws.Range("L127").Value = 4327.9092
$ws.Range("N127").Value = -14247.9092
$ws.Range("H133").Value = 25595.342
$ws.Range("I133").Value = 12782.286
$ws.Range("J133").Value = 28233.324
$ws.Range("K133").Value = 38346.858
$ws.Range("L133").Value = 84699.97200000001
$ws.Range("M133").Value = -33286.858
$ws.Range("N133").Value = -94819.97200000001
$ws.Range("H135").Value = 1080551.2
$ws.Range("I135").Value = 89500.89
$ws.Range("K135").Value = 805508.01
$ws.Range("M135").Value = -802973.01
$ws.Range("H140").Value = 2951.5
$ws.Range("I140").Value = 2368.8096
$ws.Range("J140").Value = 5398.8
$ws.Range("K140").Value = 7106.4288
$ws.Range("L140").Value = 16196.4
$ws.Range("M140").Value = -1926.4288
$ws.Range("N140").Value = -26556.4
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 309.2
$ws.Range("I2").Value = 309.2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 309.2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -196.2
$ws.Range("N2").ClearContents()
$ws.Range("H70").Value = 9833.143
$ws.Range("I70").Value = 4328.5557
$ws.Range("K70").Value = 4328.5557
$ws.Range("M70").Value = -4058.5557
$ws.Range("H73").Value = 9833.143
$ws.Range("I73").Value = 4328.5557
$ws.Range("K73").Value = 4328.5557
$ws.Range("M73").Value = -3392.5557
$ws.Range("H80").Value = 720649.5600000001
$ws.Range("I80").Value = 593869.4399999999
$ws.Range("J80").Value = 916582.5600000001
$ws.Range("K80").Value = 593869.4399999999
$ws.Range("L80").Value = 916582.5600000001
$ws.Range("M80").Value = -592871.4399999999
$ws.Range("N80").Value = -918578.5600000001
$ws.Range("H83").Value = 720649.5600000001
$ws.Range("I83").Value = 593869.4399999999
$ws.Range("J83").Value = 916582.5600000001
$ws.Range("K83").Value = 2969347.2
$ws.Range("L83").Value = 4582912.800000001
$ws.Range("M83").Value = -2964355.2
$ws.Range("N83").Value = -4592896.800000001
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H107").Value = 2475.2
$ws.Range("I107").Value = 2378.9
$ws.Range("J107").Value = 2667.8
$ws.Range("K107").Value = 2378.9
$ws.Range("L107").Value = 2667.8
$ws.Range("M107").Value = -458.9000000000001
$ws.Range("N107").Value = -6507.8
$ws.Range("H122").Value = 4717.364
$ws.Range("I122").Value = 3958.8
$ws.Range("J122").Value = 5349.5
$ws.Range("K122").Value = 11876.4
$ws.Range("L122").Value = 16048.5
$ws.Range("M122").Value = -9426.400000000001
$ws.Range("N122").Value = -20948.5
$ws.Range("H123").Value = 74989
$ws.Range("J123").Value = 74989
$ws.Range("L123").Value = 74989
$ws.Range("N123").Value = -79889
$ws.Range("H126").Value = 3418.1333
$ws.Range("I126").Value = 2549.5715
$ws.Range("J126").Value = 4178.125
$ws.Range("K126").Value = 7648.7145
$ws.Range("L126").Value = 12534.375
$ws.Range("M126").Value = -5178.7145
$ws.Range("N126").Value = -17474.375
$ws.Range("H132").Value = 60420.434
$ws.Range("I132").Value = 6860
$ws.Range("K132").Value = 20580
$ws.Range("M132").Value = -18050
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5222.3076
$ws.Range("I7").Value = 4318.8
$ws.Range("J7").Value = 5787
$ws.Range("K7").Value = 4318.8
$ws.Range("L7").Value = 5787
$ws.Range("M7").Value = -4206.8
$ws.Range("N7").Value = -6011
$ws.Range("H45").Value = 1000000
$ws.Range("J45").Value = 1000000
$ws.Range("L45").Value = 1000000
$ws.Range("N45").Value = -1000814
$ws.Range("H46").Value = 2492.1428
$ws.Range("I46").Value = 2539.9
$ws.Range("J46").Value = 2372.75
$ws.Range("K46").Value = 2539.9
$ws.Range("L46").Value = 2372.75
$ws.Range("M46").Value = -2351.9
$ws.Range("N46").Value = -2748.75
$ws.Range("H55").Value = 809.2222
$ws.Range("I55").Value = 224.4762
$ws.Range("J55").Value = 2855.8333
$ws.Range("K55").Value = 224.4762
$ws.Range("L55").Value = 2855.8333
$ws.Range("M55").Value = -51.47620000000001
$ws.Range("N55").Value = -3201.8333
$ws.Range("H61").Value = 1364.1428
$ws.Range("I61").Value = 1364.1428
$ws.Range("K61").Value = 1364.1428
$ws.Range("M61").Value = -1162.1428
$ws.Range("H82").Value = 1322.7273
$ws.Range("I82").Value = 1452.7778
$ws.Range("K82").Value = 1452.7778
$ws.Range("M82").Value = -1091.7778
$ws.Range("H85").Value = 1322.7273
$ws.Range("I85").Value = 1452.7778
$ws.Range("K85").Value = 1452.7778
$ws.Range("M85").Value = -204.7778000000001
$ws.Range("H100").Value = 3896.6
$ws.Range("I100").Value = 2494.3333
$ws.Range("K100").Value = 2494.3333
$ws.Range("M100").Value = -1953.3333
$ws.Range("H113").Value = 1364.1428
$ws.Range("I113").Value = 1364.1428
$ws.Range("K113").Value = 1364.1428
$ws.Range("M113").Value = 805.8571999999999
$ws.Range("H122").Value = 3336667.8
$ws.Range("I122").Value = 3336667.8
$ws.Range("K122").Value = 10010003.4
$ws.Range("M122").Value = -10007553.4
$ws.Range("H126").Value = 5222.3076
$ws.Range("I126").Value = 4318.8
$ws.Range("J126").Value = 5787
$ws.Range("K126").Value = 12956.4
$ws.Range("L126").Value = 17361
$ws.Range("M126").Value = -10486.4
$ws.Range("N126").Value = -22301
$ws.Range("H132").Value = 5640.067
$ws.Range("I132").Value = 4305.0625
$ws.Range("J132").Value = 7165.7856
$ws.Range("K132").Value = 12915.1875
$ws.Range("L132").Value = 21497.3568
$ws.Range("M132").Value = -10385.1875
$ws.Range("N132").Value = -26557.3568
$ws.Range("H136").Value = 922681.4399999999
$ws.Range("I136").Value = 2509999.8
$ws.Range("J136").Value = 15642.429
$ws.Range("K136").Value = 7529999.399999999
$ws.Range("L136").Value = 46927.287
$ws.Range("M136").Value = -7527449.399999999
$ws.Range("N136").Value = -52027.287
$ws.Range("H139").Value = 51199.8
$ws.Range("J139").Value = 51750
$ws.Range("L139").Value = 51750
$ws.Range("N139").Value = -62030
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 199.88889
$ws.Range("I4").Value = 157
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 157
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = -44
$ws.Range("N4").Value = -576
$ws.Range("H11").Value = 6000
$ws.Range("I11").Value = 6000
$ws.Range("K11").Value = 6000
$ws.Range("M11").Value = -5858
$ws.Range("H63").Value = 60780
$ws.Range("J63").Value = 60780
$ws.Range("L63").Value = 60780
$ws.Range("N63").Value = -62028
$ws.Range("H66").Value = 60780
$ws.Range("J66").Value = 60780
$ws.Range("L66").Value = 182340
$ws.Range("N66").Value = -188580
$ws.Range("H68").Value = 32830
$ws.Range("J68").Value = 32830
$ws.Range("L68").Value = 32830
$ws.Range("N68").Value = -34452
$ws.Range("H71").Value = 32830
$ws.Range("J71").Value = 32830
$ws.Range("L71").Value = 98490
$ws.Range("N71").Value = -106602
$ws.Range("H96").Value = 168317.33
$ws.Range("I96").Value = 334967.34
$ws.Range("K96").Value = 334967.34
$ws.Range("M96").Value = -333594.34
$ws.Range("H100").Value = 922.2
$ws.Range("I100").Value = 731.2857
$ws.Range("J100").Value = 1367.6666
$ws.Range("K100").Value = 1462.5714
$ws.Range("L100").Value = 2735.3332
$ws.Range("M100").Value = -921.5714
$ws.Range("N100").Value = -3817.3332
$ws.Range("H107").Value = 1645
$ws.Range("I107").Value = 1876.6666
$ws.Range("J107").Value = 950
$ws.Range("K107").Value = 5629.9998
$ws.Range("L107").Value = 2850
$ws.Range("M107").Value = -3709.9998
$ws.Range("N107").Value = -6690
$ws.Range("H122").Value = 18871320
$ws.Range("I122").Value = 27029678
$ws.Range("J122").Value = 5117.625
$ws.Range("K122").Value = 81089034
$ws.Range("L122").Value = 15352.875
$ws.Range("M122").Value = -81086584
$ws.Range("N122").Value = -20252.875
$ws.Range("H126").Value = 1639.6842
$ws.Range("I126").Value = 1540.875
$ws.Range("J126").Value = 2166.6667
$ws.Range("K126").Value = 4622.625
$ws.Range("L126").Value = 6500.000100000001
$ws.Range("M126").Value = -2152.625
$ws.Range("N126").Value = -11440.0001
$ws.Range("H132").Value = 143156.42
$ws.Range("I132").Value = 283.33334
$ws.Range("K132").Value = 850.0000200000001
$ws.Range("M132").Value = 1679.99998
$ws.Range("H136").Value = 310986.03
$ws.Range("I136").Value = 326050.53
$ws.Range("K136").Value = 978151.5900000001
$ws.Range("M136").Value = -975601.5900000001
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
$ws.Range("H139").Value = 57000
$ws.Range("J139").Value = 58333.332
$ws.Range("L139").Value = 58333.332
$ws.Range("N139").Value = -68613.33199999999
